# Edit script for "Linux Progress.docx":
#  1) Merge the split run in the last (Sept 22) paragraph into one run.
#  2) Append a blank paragraph, a new bold+underlined date heading
#     ("September 26th, 2022"), a new body paragraph describing the SSD/
#     theme/icons work, and a trailing blank paragraph.
#  3) Mark the section's page size as explicitly portrait.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Step 1: merge the two runs of the last paragraph into a single run ---
$lastPara = $d.Paragraphs.Last
$targetRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$mergedRunXml = '<w:p><w:r><w:rPr/><w:t>I have installed Ubuntu onto Virtual Box and resolved Guest additions issue. Must remove guest additions from storage, then power up and log into Ubuntu VM, then I must install guest additions and restart VM for it to work.</w:t></w:r></w:p>'
$targetRange.InsertXML($pkgHeader + $mergedRunXml + $pkgFooter)

# --- Step 2: append the new paragraphs ---
$blankNormalPara  = '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr></w:p>'
$dateHeadingPara  = '<w:p><w:pPr><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:u w:val="single"/></w:rPr><w:t>September 26th, 2022</w:t></w:r></w:p>'
$updatePara       = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve">I have also installed Ubuntu image onto SSD. Then, I changed the theme of Ubuntu by </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:u w:val="none"/></w:rPr><w:t>adding .themes</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:u w:val="none"/></w:rPr><w:t>and .icons</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> hidden directories to home directory and using the gnome-tweaks application. Furthermore, I went over how to navigate the file system using cd and about the ls command and its different options such as –l, -r, -p, -s.</w:t></w:r></w:p>'
$trailingBlankPara = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b w:val="0"/><w:bCs w:val="0"/><w:u w:val="none"/></w:rPr></w:pPr></w:p>'

$newBodyXml = $blankNormalPara + $dateHeadingPara + $updatePara + $trailingBlankPara

$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$insertPoint.InsertXML($pkgHeader + $newBodyXml + $pkgFooter)

# --- Step 3: mark the section page size as explicitly portrait ---
$section = $d.Sections(1)
$section.PageSetup.Orientation = 0
